$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Yesenia"
$ws.Range("B4").Value = "Gerlach"

$ws.Range("B5").Select()
